$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B2").Value = 15.90384434493724
$ws.Range("D2").Value = 3.103572750118113
$ws.Range("E2").Value = 24.05113964721559
$ws.Range("F2").Value = 18.36040925798012
$ws.Range("G2").Value = 19.59819184780265
$ws.Range("H2").Value = 11.30564026009323
$ws.Range("L2").Value = 8.083661080885756
$ws.Range("M2").Value = 14.57929188242424
$ws.Range("N2").Value = 19.35730269148822
$ws.Range("O2").Value = 16.04049802537432

$ws.Range("B3").Value = 15.74027972872933
$ws.Range("D3").Value = 3.046928811683786
$ws.Range("E3").Value = 24.01784591639422
$ws.Range("F3").Value = 18.1752103667262
$ws.Range("G3").Value = 19.1616636218117
$ws.Range("H3").Value = 11.31059940016692
$ws.Range("L3").Value = 7.987669091103283
$ws.Range("M3").Value = 14.49494885755195
$ws.Range("N3").Value = 19.36571072732071
$ws.Range("O3").Value = 15.97497286521739

$ws.Range("B4").Value = 15.64156227337502
$ws.Range("D4").Value = 3.011120286658736
$ws.Range("E4").Value = 23.99978604196047
$ws.Range("F4").Value = 18.06654879506367
$ws.Range("G4").Value = 18.89634775773411
$ws.Range("H4").Value = 11.31600835068109
$ws.Range("L4").Value = 7.928894326876406
$ws.Range("M4").Value = 14.44466226350906
$ws.Range("N4").Value = 19.3732801117641
$ws.Range("O4").Value = 15.93922921968092

$ws.Range("B5").Value = 15.60180526187795
$ws.Range("D5").Value = 2.99627903933258
$ws.Range("E5").Value = 23.99303325685025
$ws.Range("F5").Value = 18.02359142012745
$ws.Range("G5").Value = 18.78910173280559
$ws.Range("H5").Value = 11.31880718484399
$ws.Range("L5").Value = 7.905007395738019
$ws.Range("M5").Value = 14.42456244208276
$ws.Range("N5").Value = 19.37697210017775
$ws.Range("O5").Value = 15.92580431028496

$ws.Range("B6").Value = 15.59523326973685
$ws.Range("D6").Value = 2.99379991732151
$ws.Range("E6").Value = 23.99194882712109
$ws.Range("F6").Value = 18.01653977979324
$ws.Range("G6").Value = 18.77135188774542
$ws.Range("H6").Value = 11.31930784457573
$ws.Range("L6").Value = 7.901045536752958
$ws.Range("M6").Value = 14.42124901122462
$ws.Range("N6").Value = 19.37762189945694
$ws.Range("O6").Value = 15.9236443441306

$ws.Range("B7").Value = 15.64102413559859
$ws.Range("D7").Value = 3.010921127324259
$ws.Range("E7").Value = 23.99969250487188
$ws.Range("F7").Value = 18.06596403204027
$ws.Range("G7").Value = 18.89489761584381
$ws.Range("H7").Value = 11.31604368900212
$ws.Range("L7").Value = 7.928571889040463
$ws.Range("M7").Value = 14.44438958203934
$ws.Range("N7").Value = 19.37332744104917
$ws.Range("O7").Value = 15.93904353238325

$ws.Range("B8").Value = 15.8471166059098
$ws.Range("D8").Value = 3.084261174704217
$ws.Range("E8").Value = 24.03917010556239
$ws.Range("F8").Value = 18.29553692024085
$ws.Range("G8").Value = 19.44722664546512
$ws.Range("H8").Value = 11.30685970870024
$ws.Range("L8").Value = 8.050542570717068
$ws.Range("M8").Value = 14.54990821463779
$ws.Range("N8").Value = 19.35970359116248
$ws.Range("O8").Value = 16.01698100336662

$ws.Range("B9").Value = 16.26302688203869
$ws.Range("D9").Value = 3.219513335942231
$ws.Range("E9").Value = 24.13514709107869
$ws.Range("F9").Value = 18.78334674531012
$ws.Range("G9").Value = 20.54370853690504
$ws.Range("H9").Value = 11.30758921775757
$ws.Range("L9").Value = 8.290031922164227
$ws.Range("M9").Value = 14.76805025609931
$ws.Range("N9").Value = 19.35198520176308
$ws.Range("O9").Value = 16.20482648426681

$ws.Range("B10").Value = 16.57330286342478
$ws.Range("D10").Value = 3.313155632841554
$ws.Range("E10").Value = 24.21650483038011
$ws.Range("F10").Value = 19.16121779129206
$ws.Range("G10").Value = 21.34688898336983
$ws.Range("H10").Value = 11.31951020474636
$ws.Range("L10").Value = 8.464840702678739
$ws.Range("M10").Value = 14.93421552576388
$ws.Range("N10").Value = 19.35775878306458
$ws.Range("O10").Value = 16.3632587653782

$ws.Range("B11").Value = 16.71495200648786
$ws.Range("D11").Value = 3.354415500595632
$ws.Range("E11").Value = 24.25576584762604
$ws.Range("F11").Value = 19.33658363333243
$ws.Range("G11").Value = 21.70967086247607
$ws.Range("H11").Value = 11.32739150924327
$ws.Range("L11").Value = 8.543841021650884
$ws.Range("M11").Value = 15.01088047483441
$ws.Range("N11").Value = 19.36283871517678
$ws.Range("O11").Value = 16.43954148886688

$ws.Range("B12").Value = 16.76861951177456
$ws.Range("D12").Value = 3.369839612814614
$ws.Range("E12").Value = 24.27094736145874
$ws.Range("F12").Value = 19.4034238706913
$ws.Range("G12").Value = 21.84651378363229
$ws.Range("H12").Value = 11.33072776966454
$ws.Range("L12").Value = 8.573658928608396
$ws.Range("M12").Value = 15.04004816937924
$ws.Range("N12").Value = 19.36511224925992
$ws.Range("O12").Value = 16.46901223623939

$ws.Range("B13").Value = 16.75706070538549
$ws.Range("D13").Value = 3.366526770075327
$ws.Range("E13").Value = 24.2676639288702
$ws.Range("F13").Value = 19.38901042419662
$ws.Range("G13").Value = 21.81706846950559
$ws.Range("H13").Value = 11.32999363141433
$ws.Range("L13").Value = 8.567241816134935
$ws.Range("M13").Value = 15.03376061407585
$ws.Range("N13").Value = 19.36460709041437
$ws.Range("O13").Value = 16.46263955721902

$ws.Range("B14").Value = 16.7193669538065
$ws.Range("D14").Value = 3.355688503884687
$ws.Range("E14").Value = 24.25700859747806
$ws.Range("F14").Value = 19.34207432052893
$ws.Range("G14").Value = 21.72094062187344
$ws.Range("H14").Value = 11.32765894613615
$ws.Range("L14").Value = 8.546296247680948
$ws.Range("M14").Value = 15.01327749164515
$ws.Range("N14").Value = 19.36301876547827
$ws.Range("O14").Value = 16.44195448484212

$ws.Range("B15").Value = 16.6962807912639
$ws.Range("D15").Value = 3.34902347559604
$ws.Range("E15").Value = 24.2505225352648
$ws.Range("F15").Value = 19.31337902917806
$ws.Range("G15").Value = 21.6619852607039
$ws.Range("H15").Value = 11.32627463990488
$ws.Range("L15").Value = 8.533453077423419
$ws.Range("M15").Value = 15.00074819689285
$ws.Range("N15").Value = 19.36209134515464
$ws.Range("O15").Value = 16.42935969691297

$ws.Range("B16").Value = 16.5640520823714
$ws.Range("D16").Value = 3.310431605706726
$ws.Range("E16").Value = 24.21398352138808
$ws.Range("F16").Value = 19.14982141310291
$ws.Range("G16").Value = 21.32311412123898
$ws.Range("H16").Value = 11.31904447069679
$ws.Range("L16").Value = 8.459665392250713
$ws.Range("M16").Value = 14.92922537179976
$ws.Range("N16").Value = 19.35747592022611
$ws.Range("O16").Value = 16.35835644448779

$ws.Range("B17").Value = 16.48303120886884
$ws.Range("D17").Value = 3.28640833381453
$ws.Range("E17").Value = 24.19213787957298
$ws.Range("F17").Value = 19.05032605185605
$ws.Range("G17").Value = 21.1144457815001
$ws.Range("H17").Value = 11.31523750950469
$ws.Range("L17").Value = 8.414249121754244
$ws.Range("M17").Value = 14.885610932293
$ws.Range("N17").Value = 19.35527085906505
$ws.Range("O17").Value = 16.31586195558015

$ws.Range("B18").Value = 16.43648011472845
$ws.Range("D18").Value = 3.272465222663903
$ws.Range("E18").Value = 24.17978529689678
$ws.Range("F18").Value = 18.99342958238819
$ws.Range("G18").Value = 20.99419067416821
$ws.Range("H18").Value = 11.31327936347587
$ws.Range("L18").Value = 8.388079171889599
$ws.Range("M18").Value = 14.86062740445387
$ws.Range("N18").Value = 19.3542336852385
$ws.Range("O18").Value = 16.29181797792005

$ws.Range("B19").Value = 16.42072865502541
$ws.Range("D19").Value = 3.267723000564836
$ws.Range("E19").Value = 24.1756396985861
$ws.Range("F19").Value = 18.97422421057568
$ws.Range("G19").Value = 20.95343918456549
$ws.Range("H19").Value = 11.31265618016271
$ws.Range("L19").Value = 8.379210988350335
$ws.Range("M19").Value = 14.85218654402849
$ws.Range("N19").Value = 19.35392229826014
$ws.Range("O19").Value = 16.28374604330234

$ws.Range("B20").Value = 16.49165117016836
$ws.Range("D20").Value = 3.288978704058494
$ws.Range("E20").Value = 24.19444146530543
$ws.Range("F20").Value = 19.06088376176084
$ws.Range("G20").Value = 21.13668437146926
$ws.Range("H20").Value = 11.31561881744486
$ws.Range("L20").Value = 8.419088866870418
$ws.Range("M20").Value = 14.89024330370871
$ws.Range("N20").Value = 19.35548169249602
$ws.Range("O20").Value = 16.32034455440326

$ws.Range("B21").Value = 16.73043811103315
$ws.Range("D21").Value = 3.358877454459615
$ws.Range("E21").Value = 24.26012987529717
$ws.Range("F21").Value = 19.35584936720245
$ws.Range("G21").Value = 21.74919145586733
$ws.Range("H21").Value = 11.32833516780637
$ws.Range("L21").Value = 8.552451302355808
$ws.Range("M21").Value = 15.01929032775934
$ws.Range("N21").Value = 19.3634758235787
$ws.Range("O21").Value = 16.44801451517582

$ws.Range("B22").Value = 16.88664260204531
$ws.Range("D22").Value = 3.403390785132123
$ws.Range("E22").Value = 24.30488801721097
$ws.Range("F22").Value = 19.55111969340626
$ws.Range("G22").Value = 22.14631560623758
$ws.Range("H22").Value = 11.33869547876124
$ws.Range("L22").Value = 8.639031151307828
$ws.Range("M22").Value = 15.10441612200784
$ws.Range("N22").Value = 19.37073833019564
$ws.Range("O22").Value = 16.53484827494316

$ws.Range("B23").Value = 16.8032746662462
$ws.Range("D23").Value = 3.379742586351906
$ws.Range("E23").Value = 24.28083576865564
$ws.Range("F23").Value = 19.44669393127981
$ws.Range("G23").Value = 21.93470534362979
$ws.Range("H23").Value = 11.33297910988296
$ws.Range("L23").Value = 8.592882340612476
$ws.Range("M23").Value = 15.05891704727192
$ws.Range("N23").Value = 19.36667673124533
$ws.Range("O23").Value = 16.48820027592404

$ws.Range("B24").Value = 16.48775399228936
$ws.Range("D24").Value = 3.287817049653147
$ws.Range("E24").Value = 24.19339936944375
$ws.Range("F24").Value = 19.05610966683134
$ws.Range("G24").Value = 21.12663119228896
$ws.Range("H24").Value = 11.31544570990978
$ws.Range("L24").Value = 8.416901001776724
$ws.Range("M24").Value = 14.88814872374909
$ws.Range("N24").Value = 19.35538565635544
$ws.Range("O24").Value = 16.31831676539738

$ws.Range("B25").Value = 16.14950651031204
$ws.Range("D25").Value = 3.183895273070362
$ws.Range("E25").Value = 24.1072434077814
$ws.Range("F25").Value = 18.64770387971093
$ws.Range("G25").Value = 20.24676703261735
$ws.Range("H25").Value = 11.30538860066325
$ws.Range("L25").Value = 8.225349361575757
$ws.Range("M25").Value = 14.70793173738156
$ws.Range("N25").Value = 19.35205305502728
$ws.Range("O25").Value = 16.15035033855728
